$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Responsibile Parties")
$ws.Activate()

# Insert a new row before row 10, shifting existing rows (10-34) down to (11-35)
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with Doroteaciro Iovino's details
$ws.Cells.Item(10, 1).Value = "IOVINO-DOROTEACIRO"
$ws.Cells.Item(10, 2).Value = "Doroteaciro Iovino "
$ws.Cells.Item(10, 3).Value = $false
$ws.Cells.Item(10, 4).Value = "Viale Carlo Berti Pichat, 6/2`n40127, Bologna - Italy"
$ws.Cells.Item(10, 5).Value = "dorotea.iovino@cmcc.it"
$ws.Cells.Item(10, 6).Value = "https://www.cmcc.it/people/iovino-doroteaciro"
$ws.Cells.Item(10, 7).Value = "https://orcid.org/0000-0001-5132-7255"

$ws.Range("A10").Font.Bold = $true
$ws.Range("A10:B10").HorizontalAlignment = -4131

# Update active cell selection to A13 (as in the target workbook)
$ws.Range("A13").Select()
